$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a plain decimal number (e.g. "0.541") need an
# explicit Text number format first, otherwise Excel auto-converts the assigned
# string into a numeric value instead of keeping it as text.

$ws.Cells.Item(2, 4).Value = '25.804.78'
$ws.Cells.Item(2, 5).Value = '  +0.24%  '

$ws.Cells.Item(3, 4).Value = '1.623.63'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.22%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '214.24'
$ws.Cells.Item(5, 5).Value = '  -0.50%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.499'
$ws.Cells.Item(6, 5).Value = '  -0.11%  '

$ws.Cells.Item(7, 5).Value = '  -0.25%  '

$ws.Cells.Item(8, 5).Value = '  -0.54%  '

$ws.Cells.Item(9, 5).Value = '  -0.38%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.57'
$ws.Cells.Item(10, 5).Value = '  +0.56%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0786'
$ws.Cells.Item(11, 5).Value = '  -0.78%  '

$ws.Cells.Item(12, 4).Value = '1.847.54'
$ws.Cells.Item(12, 5).Value = '  -0.55%  '

$ws.Cells.Item(13, 5).Value = '  -0.66%  '

$ws.Cells.Item(14, 4).Value = '1.620.14'
$ws.Cells.Item(14, 5).Value = '  -0.46%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.541'
$ws.Cells.Item(15, 5).Value = '  -2.92%  '

$ws.Cells.Item(16, 4).Value = '0.0₃0755'
$ws.Cells.Item(16, 5).Value = '  -0.92%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '62.44'
$ws.Cells.Item(17, 5).Value = '  -0.95%  '

$ws.Cells.Item(18, 4).Value = '25.787.08'
$ws.Cells.Item(18, 5).Value = '  +0.09%  '

$ws.Cells.Item(19, 5).Value = '  -0.21%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '192.17'
$ws.Cells.Item(20, 5).Value = '  +0.07%  '

$ws.Cells.Item(21, 5).Value = '  -2.21%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.92'
$ws.Cells.Item(22, 5).Value = '  -0.25%  '

$ws.Cells.Item(23, 5).Value = '  -0.80%  '

$ws.Cells.Item(24, 5).Value = '  -0.97%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.998'
$ws.Cells.Item(25, 5).Value = '  -0.34%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '141.68'
$ws.Cells.Item(26, 5).Value = '  -0.91%  '

$ws.Cells.Item(27, 5).Value = '  +1.27%  '

$ws.Cells.Item(28, 5).Value = '  -0.37%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '15.41'
$ws.Cells.Item(29, 5).Value = '  -0.40%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.23'
$ws.Cells.Item(30, 5).Value = '  -0.46%  '

$ws.Cells.Item(31, 5).Value = '  +0.87%  '

$ws.Cells.Item(32, 5).Value = '  -0.75%  '

$ws.Cells.Item(33, 5).Value = '  -0.66%  '

$ws.Cells.Item(34, 5).Value = '  +0.22%  '

$ws.Cells.Item(35, 5).Value = '  +1.21%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.899'
$ws.Cells.Item(36, 5).Value = '  -0.37%  '

$ws.Cells.Item(37, 4).Value = '1.124.62'
$ws.Cells.Item(37, 5).Value = '  -0.55%  '

$ws.Cells.Item(38, 5).Value = '  +0.29%  '

$ws.Cells.Item(39, 5).Value = '  -2.13%  '

$ws.Cells.Item(40, 5).Value = '  +0.60%  '

$ws.Cells.Item(41, 5).Value = '  -0.31%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '99.24'
$ws.Cells.Item(42, 5).Value = '  -1.43%  '

$ws.Cells.Item(43, 5).Value = '  -2.35%  '

$ws.Cells.Item(44, 5).Value = '  -0.43%  '

$ws.Cells.Item(45, 4).Value = '1.758.63'
$ws.Cells.Item(45, 5).Value = '  -0.46%  '

$ws.Cells.Item(46, 4).Value = '0.0₆0110'
$ws.Cells.Item(46, 5).Value = '  -1.42%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '56.14'
$ws.Cells.Item(47, 5).Value = '  +1.47%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0524'
$ws.Cells.Item(48, 5).Value = '  +3.21%  '

$ws.Cells.Item(49, 5).Value = '  +2.99%  '

$ws.Cells.Item(50, 5).Value = '  -0.75%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.57'
$ws.Cells.Item(51, 5).Value = '  +1.66%  '
